# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G holds the "K" (strikeouts) stat. The underlying game log was
# regenerated from source (K replacing the old Strike# derived value), so
# the freshly calculated s_vals for K are written back into the sheet,
# row by row, for every game record (rows 2-69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly calculated K values (s_vals), in row order for rows 2..69.
$kValues = @(
    2, 0, 2, 1, 1, 1, 1, 2, 1, 0,
    1, 1, 2, 1, 1, 2, 0, 1, 1, 0,
    0, 0, 2, 0, 1, 1, 1, 1, 1, 0,
    0, 0, 1, 0, 1, 0, 1, 1, 2, 0,
    2, 0, 1, 0, 2, 3, 2, 1, 1, 2,
    2, 1, 1, 1, 1, 2, 2, 1, 1, 1,
    1, 1, 1, 0, 0, 0, 0, 1
)

$firstRow = 2
$lastRow = 69
$rowCount = $kValues.Length

$kRange = New-Object 'object[,]' $rowCount, 1
for ($i = 0; $i -lt $rowCount; $i++) {
    $kRange[$i, 0] = $kValues[$i]
}

$ws.Range("G$firstRow`:G$lastRow").Value = $kRange
